$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-05-28 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-29 Monday", 2) | Out-Null

# Update table cell values by position (row, col), 1-based
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "6-1=5"
$t.Cell(1, 2).Range.Text = "99-17=82"
$t.Cell(1, 3).Range.Text = "81-10=71"
$t.Cell(1, 4).Range.Text = "9+16=25"
$t.Cell(1, 5).Range.Text = "14+8=22"
$t.Cell(2, 1).Range.Text = "59-45=14"
$t.Cell(2, 2).Range.Text = "51+26=77"
$t.Cell(2, 3).Range.Text = "53-14=39"
$t.Cell(2, 4).Range.Text = "64+31=95"
$t.Cell(2, 5).Range.Text = "15+34=49"
$t.Cell(3, 1).Range.Text = "54-10=44"
$t.Cell(3, 2).Range.Text = "69-13=56"
$t.Cell(3, 3).Range.Text = "52-25=27"
$t.Cell(3, 4).Range.Text = "80+12=92"
$t.Cell(3, 5).Range.Text = "89-10=79"
$t.Cell(4, 1).Range.Text = "14+51=65"
$t.Cell(4, 2).Range.Text = "18-2=16"
$t.Cell(4, 3).Range.Text = "80+19=99"
$t.Cell(4, 4).Range.Text = "98-37=61"
$t.Cell(4, 5).Range.Text = "34-27=7"
$t.Cell(5, 1).Range.Text = "10+46=56"
$t.Cell(5, 2).Range.Text = "46+26=72"
$t.Cell(5, 3).Range.Text = "76+12=88"
$t.Cell(5, 4).Range.Text = "89-65=24"
$t.Cell(5, 5).Range.Text = "15+50=65"
$t.Cell(6, 1).Range.Text = "55-36=19"
$t.Cell(6, 2).Range.Text = "95-77=18"
$t.Cell(6, 3).Range.Text = "45-27=18"
$t.Cell(6, 4).Range.Text = "83-66=17"
$t.Cell(6, 5).Range.Text = "89-47=42"
$t.Cell(7, 1).Range.Text = "27-8=19"
$t.Cell(7, 2).Range.Text = "95-14=81"
$t.Cell(7, 3).Range.Text = "60-7=53"
$t.Cell(7, 4).Range.Text = "83-46=37"
$t.Cell(7, 5).Range.Text = "8+60=68"
$t.Cell(8, 1).Range.Text = "46+20=66"
$t.Cell(8, 2).Range.Text = "11+54=65"
$t.Cell(8, 3).Range.Text = "50+49=99"
$t.Cell(8, 4).Range.Text = "33-26=7"
$t.Cell(8, 5).Range.Text = "95-8=87"
$t.Cell(9, 1).Range.Text = "19+8=27"
$t.Cell(9, 2).Range.Text = "51-27=24"
$t.Cell(9, 3).Range.Text = "16+55=71"
$t.Cell(9, 4).Range.Text = "88+4=92"
$t.Cell(9, 5).Range.Text = "7+0=7"
$t.Cell(10, 1).Range.Text = "94-25=69"
$t.Cell(10, 2).Range.Text = "75-43=32"
$t.Cell(10, 3).Range.Text = "88+0=88"
$t.Cell(10, 4).Range.Text = "73+7=80"
$t.Cell(10, 5).Range.Text = "62-37=25"
$t.Cell(11, 1).Range.Text = "16+53=69"
$t.Cell(11, 2).Range.Text = "90-0=90"
$t.Cell(11, 3).Range.Text = "19+42=61"
$t.Cell(11, 4).Range.Text = "24+55=79"
$t.Cell(11, 5).Range.Text = "20+42=62"
$t.Cell(12, 1).Range.Text = "18+2=20"
$t.Cell(12, 2).Range.Text = "19+18=37"
$t.Cell(12, 3).Range.Text = "62-11=51"
$t.Cell(12, 4).Range.Text = "68-5=63"
$t.Cell(12, 5).Range.Text = "44+20=64"
$t.Cell(13, 1).Range.Text = "75+21=96"
$t.Cell(13, 2).Range.Text = "69-48=21"
$t.Cell(13, 3).Range.Text = "49+9=58"
$t.Cell(13, 4).Range.Text = "36+52=88"
$t.Cell(13, 5).Range.Text = "30-23=7"
$t.Cell(14, 1).Range.Text = "61-40=21"
$t.Cell(14, 2).Range.Text = "16+24=40"
$t.Cell(14, 3).Range.Text = "88-86=2"
$t.Cell(14, 4).Range.Text = "3+13=16"
$t.Cell(14, 5).Range.Text = "50-11=39"
$t.Cell(15, 1).Range.Text = "67-11=56"
$t.Cell(15, 2).Range.Text = "57+4=61"
$t.Cell(15, 3).Range.Text = "3+70=73"
$t.Cell(15, 4).Range.Text = "73-53=20"
$t.Cell(15, 5).Range.Text = "32+33=65"
$t.Cell(16, 1).Range.Text = "36+49=85"
$t.Cell(16, 2).Range.Text = "85-61=24"
$t.Cell(16, 3).Range.Text = "50+33=83"
$t.Cell(16, 4).Range.Text = "47+37=84"
$t.Cell(16, 5).Range.Text = "45-26=19"
$t.Cell(17, 1).Range.Text = "27-10=17"
$t.Cell(17, 2).Range.Text = "60+16=76"
$t.Cell(17, 3).Range.Text = "85-20=65"
$t.Cell(17, 4).Range.Text = "95-24=71"
$t.Cell(17, 5).Range.Text = "97-90=7"
$t.Cell(18, 1).Range.Text = "71-10=61"
$t.Cell(18, 2).Range.Text = "27+42=69"
$t.Cell(18, 3).Range.Text = "50+8=58"
$t.Cell(18, 4).Range.Text = "65+4=69"
$t.Cell(18, 5).Range.Text = "65-16=49"
$t.Cell(19, 1).Range.Text = "92-90=2"
$t.Cell(19, 2).Range.Text = "51+23=74"
$t.Cell(19, 3).Range.Text = "75-22=53"
$t.Cell(19, 4).Range.Text = "21+31=52"
$t.Cell(19, 5).Range.Text = "25+56=81"
$t.Cell(20, 1).Range.Text = "2+74=76"
$t.Cell(20, 2).Range.Text = "74+1=75"
$t.Cell(20, 3).Range.Text = "18-4=14"
$t.Cell(20, 4).Range.Text = "76+13=89"
$t.Cell(20, 5).Range.Text = "5+85=90"

Write-Output "done"
